# "Adding latest code for 'On This Page' section"
#
# - Sheet "pages_with_on_this_page": drop the extra PDQ/Spanish rows (rows
#   3-6), keep just the header + the single English "about-cancer/coping/
#   self-image" article row, rename the "header" column to "otpHeader", and
#   tighten column A's width (drop the stale bestFit autosize).
# - Sheet "pages_without_on_this_page": stays the same one data row, just
#   becomes the non-selected tab with its selection reset to A2.
# - The workbook opens on the first sheet (pages_with_on_this_page) with
#   B14 selected.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("pages_with_on_this_page")
$ws2 = $wb.Worksheets.Item("pages_without_on_this_page")

# --- pages_with_on_this_page ---------------------------------------------

# Header row: "header" -> "otpHeader"
$ws1.Cells.Item(1, 4).Value = "otpHeader"

# Drop rows 3-6 (Spanish / PDQ Cancer rows), leaving header + row 2
[void]$ws1.Rows("3:6").Delete()

# Column A no longer needs the old bestFit 69.16 width; new fixed width 57
$ws1.Columns.Item(1).ColumnWidth = 56.16666666666667

# --- pages_without_on_this_page ------------------------------------------

# Reset its selection to A2 (it loses the active tab to sheet 1 below)
[void]$ws2.Range("A2").Select()

# --- Make pages_with_on_this_page the active tab/selection ---------------

[void]$ws1.Activate()
[void]$ws1.Range("B14").Select()
